$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing
# columns A:H to B:I (preserving their values/styles/types).
$ws.Columns("A:A").Insert()

# Fill the new column A with row labels.
$ws.Range("A1").Value = "index"
$ws.Range("A2").Value = "Hız"
$ws.Range("A3").Value = "Aşım Oranı"
$ws.Range("A4").Value = "Para Cezası"
$ws.Range("A5").Value = "Aşma Sınıfı"

# Give the new label cells the same formatting (bold, bordered,
# centered/top-aligned) already used by the header row, by copying
# the format from B1 (the shifted former A1 header cell).
$ws.Range("B1").Copy()
$ws.Range("A1:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
